# Update the three-digit x one-digit multiplication answers in the table.
# Each source text is unique in the document, so a simple Find/Replace
# (restricted to replace only the first/next match) per pair is safe and
# keeps run formatting (font, size) untouched.

$d = $word.ActiveDocument

$pairs = @(
    @{ old = "793×8=6344"; new = "392×7=2744" },
    @{ old = "934×2=1868"; new = "233×5=1165" },
    @{ old = "712×7=4984"; new = "978×2=1956" },
    @{ old = "859×5=4295"; new = "104×4=416" },
    @{ old = "738×6=4428"; new = "961×7=6727" },
    @{ old = "641×2=1282"; new = "890×4=3560" },
    @{ old = "728×8=5824"; new = "874×5=4370" },
    @{ old = "177×5=885";  new = "280×4=1120" },
    @{ old = "768×7=5376"; new = "729×3=2187" },
    @{ old = "172×3=516";  new = "121×4=484" },
    @{ old = "485×2=970";  new = "927×9=8343" },
    @{ old = "825×6=4950"; new = "158×7=1106" },
    @{ old = "628×8=5024"; new = "302×9=2718" },
    @{ old = "904×2=1808"; new = "588×4=2352" },
    @{ old = "575×7=4025"; new = "621×5=3105" },
    @{ old = "449×8=3592"; new = "836×6=5016" },
    @{ old = "622×2=1244"; new = "233×5=1165" },
    @{ old = "990×2=1980"; new = "316×4=1264" },
    @{ old = "173×6=1038"; new = "873×9=7857" },
    @{ old = "142×3=426";  new = "101×7=707" },
    @{ old = "254×4=1016"; new = "768×7=5376" },
    @{ old = "557×3=1671"; new = "964×3=2892" },
    @{ old = "674×4=2696"; new = "410×4=1640" },
    @{ old = "614×7=4298"; new = "446×5=2230" },
    @{ old = "478×2=956";  new = "501×7=3507" }
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $pair.new, 2)
}
